$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (blank) column before column N, shifting the existing
# "Late" / "Outstanding" columns one place to the right.
$mColWidth = $ws.Columns("M:M").ColumnWidth
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $mColWidth

# Make "Repayment schedule" the active sheet/tab and set its selection,
# matching the updated workbook view state.
$ws.Activate() | Out-Null
$ws.Range("S8").Select() | Out-Null
